$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly report data between row 2 and row 3 (dates, volumes,
# prices, units and origin measures) so each row reflects the correct
# week's figures.

$cols = @("D", "M", "N", "O", "P", "Q", "S", "T")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $tmp = $cell2.Value2
    $cell2.Value = $cell3.Value2
    $cell3.Value = $tmp
}
